$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.394.97'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').Value = '1.802.33'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '227.88'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').Value = '0.580'
$ws.Range('E6').Value = '  +3.98%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '34.97'
$ws.Range('E8').Value = '  +6.31%  '
$ws.Range('D9').Value = '0.299'
$ws.Range('E9').Value = '  +0.48%  '
$ws.Range('D10').Value = '0.0691'
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').Value = '2.063.66'
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('E13').Value = '  +0.37%  '
$ws.Range('D14').Value = '1.800.53'
$ws.Range('E14').Value = '  +0.02%  '
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('D16').Value = '34.387.54'
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('D17').Value = '4.33'
$ws.Range('E17').Value = '  +1.20%  '
$ws.Range('D18').Value = '68.90'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0796'
$ws.Range('E19').Value = '  -0.73%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '245.04'
$ws.Range('E20').Value = '  -1.32%  '
$ws.Range('D21').Value = '11.49'
$ws.Range('E21').Value = '  +1.51%  '
$ws.Range('E22').Value = '  +0.25%  '
$ws.Range('E23').Value = '  -0.77%  '
$ws.Range('D24').Value = '169.80'
$ws.Range('E24').Value = '  +2.88%  '
$ws.Range('E25').Value = '  +2.78%  '
$ws.Range('D26').Value = '7.60'
$ws.Range('E26').Value = '  +4.62%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = '0.119'
$ws.Range('E27').Value = '  +2.18%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '16.69'
$ws.Range('E28').Value = '  +0.68%  '
$ws.Range('E29').Value = '  -0.88%  '
$ws.Range('D30').Value = '3.97'
$ws.Range('E30').Value = '  -5.51%  '
$ws.Range('D31').Value = '0.0527'
$ws.Range('E31').Value = '  +0.84%  '
$ws.Range('D32').Value = '1.25'
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').Value = '1.394.64'
$ws.Range('E35').Value = '  -2.33%  '
$ws.Range('D36').Value = '0.677'
$ws.Range('E36').Value = '  +0.64%  '
$ws.Range('D37').Value = '2.52'
$ws.Range('E37').Value = '  -2.36%  '
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('E39').Value = '  -1.18%  '
$ws.Range('D40').Value = '82.79'
$ws.Range('E40').Value = '  -3.23%  '
$ws.Range('D41').Value = '2.83'
$ws.Range('E41').Value = '  +3.11%  '
$ws.Range('D42').Value = '0.944'
$ws.Range('E42').Value = '  +1.10%  '
$ws.Range('D44').Value = '13.53'
$ws.Range('E44').Value = '  -0.83%  '
$ws.Range('E45').Value = '  +3.06%  '
$ws.Range('E46').Value = '  -3.09%  '
$ws.Range('E47').Value = '  -2.01%  '
$ws.Range('D48').Value = '1.963.76'
$ws.Range('E48').Value = '  +0.44%  '
$ws.Range('D49').Value = '104.43'
$ws.Range('E49').Value = '  -1.65%  '
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('E51').Value = '  +1.01%  '
